# Apply audit spreadsheet corrections provided by prof Giselle.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5-6: Computer Science requirement gains "Introduction to Computer Systems" subcategory
$ws.Range("A5").Value = "BS in Computer Science---Computer Science---Introduction to Computer Systems"
$ws.Range("A6").Value = "BS in Computer Science---Computer Science---Introduction to Computer Systems"

# Rows 14-17: AI elective course list shifted down by one, with a new first entry
$ws.Range("B14").Value = "15-382"
$ws.Range("B15").Value = "15-386"
$ws.Range("B16").Value = "16-384"
$ws.Range("B17").Value = "16-385"

# Rows 45-48: fix capitalization of "3D Calculus"
$ws.Range("A45").Value = "BS in Computer Science---Mathematics and Probability---Calculus---3D Calculus"
$ws.Range("A46").Value = "BS in Computer Science---Mathematics and Probability---Calculus---3D Calculus"
$ws.Range("A47").Value = "BS in Computer Science---Mathematics and Probability---Calculus---3D Calculus"
$ws.Range("A48").Value = "BS in Computer Science---Mathematics and Probability---Calculus---3D Calculus"

# Rows 57-60: simplify probability & statistics sequence labels
$ws.Range("A57").Value = "BS in Computer Science---Mathematics and Probability---Probability---Probability and Statistics 36-22x sequence"
$ws.Range("A58").Value = "BS in Computer Science---Mathematics and Probability---Probability---Probability and Statistics 36-22x sequence"
$ws.Range("A59").Value = "BS in Computer Science---Mathematics and Probability---Probability---Probability and Statistics 36-23x sequence"
$ws.Range("A60").Value = "BS in Computer Science---Mathematics and Probability---Probability---Probability and Statistics 36-23x sequence"

# Rows 64-66: rename "2 SCS Electives" requirement to "SCS Electives"
$ws.Range("A64").Value = "BS in Computer Science---SCS Electives"
$ws.Range("A65").Value = "BS in Computer Science---SCS Electives"
$ws.Range("A66").Value = "BS in Computer Science---SCS Electives"
